# Apply "Updated with Rapise 6.6 note" edit:
#  1. On the "RVL" sheet, rename several "Functions"/"Nav*" rows to use the
#     shorter "Nav" object name together with the corresponding shortened
#     action name (the NAV prefix is dropped from the action column).
#  2. On the "Cleanup" sheet, remove the obsolete 4-row "NavVersion" block
#     (rows 15-18), shifting everything below it up by 4 rows.

$wb = $excel.ActiveWorkbook

# --- 1. RVL sheet: update Object/Action columns -----------------------
$rvl = $wb.Worksheets.Item("RVL")

$rvl.Range("C21").Value = "Nav"
$rvl.Range("D21").Value = "Launch"

$rvl.Range("C22").Value = "Nav"
$rvl.Range("D22").Value = "ChangeCompany"

$rvl.Range("C23").Value = "Nav"
$rvl.Range("D23").Value = "Navigate"

$rvl.Range("C37").Value = "Nav"
$rvl.Range("D37").Value = "Version"

$rvl.Range("C41").Value = "Nav"
$rvl.Range("D41").Value = "SelectFastTab"

$rvl.Range("C46").Value = "Nav"
$rvl.Range("D46").Value = "SelectFastTab"

$rvl.Range("C53").Value = "Nav"
$rvl.Range("D53").Value = "Close"

# --- 2. Cleanup sheet: delete the old NavVersion rows (15-18) ---------
$cleanup = $wb.Worksheets.Item("Cleanup")
$cleanup.Rows("15:18").Delete()
